# comprehension_questions.xlsx edits:
#  1. "Az első feladat nagyjából 20 percet..." -> "...25 percet..." (cells F4 & F5,
#     which share the same text) and bump their row height to fit the now-taller
#     wrapped text.
#  2. Collapse the two-run rich text in A12 ("Hogy érzi, ... feladatokat? " +
#     "Kérjük, ... nem befolyásolja.") into a single, uniformly formatted string.
#  3. Move the active selection to F5 (scrolled so column C is the first visible
#     column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$repeatedInfo = "Az első feladat nagyjából 25 percet vesz igénybe, közben két rövid szünettel. Ügyeljen, hogy ezek a szünetek ne legyenek 2 percnél hosszabbak. `n`nA 'J' billentyűvel jelölje azokat a képeket, amelyek maradhatnak a galériában, a bemutatott helyen.`n`nAz 'F' billentyűvel jelölje a képeket, amelyek nem maradnak kiállítva a bemutatott helyen. "

$ws.Range("F4").Value = $repeatedInfo
$ws.Range("F5").Value = $repeatedInfo

$ws.Rows.Item(4).RowHeight = 188.3
$ws.Rows.Item(5).RowHeight = 188.3

$ws.Range("A12").Value = "Hogy érzi, lelkiismeretesen, figyelmesen oldotta meg a feladatokat? Kérjük, válaszoljon őszintén. Válasza a vizsgálat bejefezésével járó jutalom (kredit, ajándékutalvány) értékét nem befolyásolja."

$ws.Range("F5").Select()
